# "Update Roles and SNMP page"
#
# 1) Global default font: Calibri -> Arial
# 2) TestData sheet: minor column-width tweaks (cols C & D)
# 3) Role sheet: minor column-width tweak (col A)
# 4) New "SNMP" worksheet appended after "Role", populated with
#    securityName/authPro/khang/MD5 and made the active sheet/cell.

$wb = $excel.ActiveWorkbook

# --- 1) Workbook-wide default font -----------------------------------
$wb.Styles("Normal").Font.Name = "Arial"

# --- 2) TestData (sheet 1) column widths -------------------------------
$wsTestData = $wb.Worksheets.Item("TestData")
$wsTestData.Columns.Item(3).ColumnWidth = 10.0    # col C
$wsTestData.Columns.Item(4).ColumnWidth = 14.25   # col D
# Default row height follows the new standard font (Calibri 11 -> Arial 10).
$wsTestData.StandardHeight = 14.25

# --- 3) Role (sheet 2) column width ------------------------------------
$wsRole = $wb.Worksheets.Item("Role")
$wsRole.Columns.Item(1).ColumnWidth = 10.0        # col A
$wsRole.StandardHeight = 14.25

# --- 4) Add the new SNMP worksheet at the end --------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsSnmp = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsSnmp.Name = "SNMP"

# Fill column-by-column (A then B) so new shared strings are interned
# in the same order as the target workbook: securityName, khang,
# authPro, MD5.
$wsSnmp.Range("A1").Value = "securityName"
$wsSnmp.Range("A2").Value = "khang"
$wsSnmp.Range("B1").Value = "authPro"
$wsSnmp.Range("B2").Value = "MD5"

$wsSnmp.StandardHeight = 14.25

# Match the recorded selection on the new sheet.
$null = $wsSnmp.Range("W8").Select()
